$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill F2:F131 with the scaling value 61 (scale by individual PPV)
$ws.Range("F2:F131").Value = 61

# Update the view state: scroll so row 123 is at the top and select C128
$excel.ActiveWindow.ScrollRow = 123
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C128").Select()
